$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2084.261
$ws.Range("J17").Value = 2084.261
$ws.Range("L17").Value = 6252.782999999999
$ws.Range("N17").Value = -6588.782999999999
# Row 18
$ws.Range("H18").Value = 774.75
$ws.Range("I18").Value = 774.75
$ws.Range("K18").Value = 774.75
$ws.Range("M18").Value = -490.75
# Row 62
$ws.Range("H62").Value = 5414
$ws.Range("I62").Value = 4691
$ws.Range("J62").Value = 5896
$ws.Range("K62").Value = 4691
$ws.Range("L62").Value = 5896
$ws.Range("M62").Value = -4067
$ws.Range("N62").Value = -7144
# Row 65
$ws.Range("H65").Value = 5414
$ws.Range("I65").Value = 4691
$ws.Range("J65").Value = 5896
$ws.Range("K65").Value = 23455
$ws.Range("L65").Value = 29480
$ws.Range("M65").Value = -20335
$ws.Range("N65").Value = -35720
# Row 106
$ws.Range("H106").Value = 4749.25
$ws.Range("I106").Value = 4666
$ws.Range("K106").Value = 4666
$ws.Range("M106").Value = -4035
# Row 137
$ws.Range("H137").Value = 3585.625
$ws.Range("I137").Value = 2564.6667
$ws.Range("K137").Value = 7694.000100000001
$ws.Range("M137").Value = -5144.000100000001
# Row 138
$ws.Range("H138").Value = 2896
$ws.Range("I138").Value = 1728.6666
$ws.Range("J138").Value = 3165.3845
$ws.Range("K138").Value = 5185.9998
$ws.Range("L138").Value = 9496.1535
$ws.Range("M138").Value = -45.9997999999996
$ws.Range("N138").Value = -19776.1535

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1075
$ws.Range("I2").Value = 602.25
$ws.Range("K2").Value = 602.25
$ws.Range("M2").Value = -489.25
# Row 33
$ws.Range("H33").Value = 11570
$ws.Range("I33").Value = 10998.333
$ws.Range("J33").Value = 15000
$ws.Range("K33").Value = 10998.333
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = -10669.333
$ws.Range("N33").Value = -15658
# Row 105
$ws.Range("H105").Value = 55554.5
$ws.Range("J105").Value = 55554.5
$ws.Range("L105").Value = 55554.5
$ws.Range("N105").Value = -62542.5
# Row 114
$ws.Range("H114").Value = 69995
$ws.Range("J114").Value = 69995
$ws.Range("L114").Value = 69995
$ws.Range("N114").Value = -78673
# Row 116
$ws.Range("H116").Value = 1075
$ws.Range("I116").Value = 602.25
$ws.Range("K116").Value = 602.25
$ws.Range("M116").Value = 1691.75
# Row 122
$ws.Range("H122").Value = 1150
$ws.Range("I122").Value = 1150
$ws.Range("K122").Value = 3450
$ws.Range("M122").Value = -1000

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1075
$ws.Range("I3").Value = 602.25
$ws.Range("K3").Value = 602.25
$ws.Range("M3").Value = -488.25
# Row 86
$ws.Range("H86").Value = 3500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 3500
$ws.Range("N86").Value = -5746
# Row 89
$ws.Range("H89").Value = 3500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 17500
$ws.Range("N89").Value = -28732
# Row 107
$ws.Range("H107").Value = 1050
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
# Row 137
$ws.Range("H137").Value = 46999.2
$ws.Range("I137").Value = 35000
$ws.Range("J137").Value = 49999
$ws.Range("K137").Value = 35000
$ws.Range("L137").Value = 49999
$ws.Range("M137").Value = -29900
$ws.Range("N137").Value = -60199

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3017.9565
$ws.Range("I31").Value = 2571.0667
$ws.Range("J31").Value = 3855.875
$ws.Range("K31").Value = 2571.0667
$ws.Range("L31").Value = 3855.875
$ws.Range("M31").Value = -2276.0667
$ws.Range("N31").Value = -4445.875
# Row 34
$ws.Range("H34").Value = 3017.9565
$ws.Range("I34").Value = 2571.0667
$ws.Range("J34").Value = 3855.875
$ws.Range("K34").Value = 2571.0667
$ws.Range("L34").Value = 3855.875
$ws.Range("M34").Value = -2369.0667
$ws.Range("N34").Value = -4259.875
# Row 50
$ws.Range("H50").Value = 20291.5
$ws.Range("I50").Value = 10583
$ws.Range("K50").Value = 10583
$ws.Range("M50").Value = -9958
# Row 58
$ws.Range("H58").Value = 1688.8125
$ws.Range("I58").Value = 1634.8
$ws.Range("K58").Value = 1634.8
$ws.Range("M58").Value = -1431.8
# Row 59
$ws.Range("H59").Value = 32825
$ws.Range("I59").Value = 28475
$ws.Range("K59").Value = 28475
$ws.Range("M59").Value = -27330
# Row 105
$ws.Range("H105").Value = 2176.2856
$ws.Range("I105").Value = 1416.6666
$ws.Range("K105").Value = 1416.6666
$ws.Range("M105").Value = 330.3334
# Row 122
$ws.Range("H122").Value = 715.2222
$ws.Range("I122").Value = 693.8570999999999
$ws.Range("K122").Value = 2081.5713
$ws.Range("M122").Value = 368.4287000000004
# Row 132
$ws.Range("H132").Value = 988.3125
$ws.Range("I132").Value = 994.6
$ws.Range("K132").Value = 2983.8
$ws.Range("M132").Value = -453.8000000000002
# Row 136
$ws.Range("H136").Value = 1688.8125
$ws.Range("I136").Value = 1634.8
$ws.Range("K136").Value = 4904.4
$ws.Range("M136").Value = -2354.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1723.7693
$ws.Range("I4").Value = 1676.4117
$ws.Range("K4").Value = 5029.2351
$ws.Range("M4").Value = -4917.2351
# Row 29
$ws.Range("H29").Value = 24.44186
$ws.Range("I29").Value = 275.5
$ws.Range("J29").Value = 12.195122
$ws.Range("K29").Value = 826.5
$ws.Range("L29").Value = 36.585366
$ws.Range("M29").Value = -549.5
$ws.Range("N29").Value = -590.585366
# Row 39
$ws.Range("H39").Value = 4693.3335
$ws.Range("J39").Value = 4693.3335
$ws.Range("L39").Value = 14080.0005
$ws.Range("N39").Value = -14668.0005
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
# Row 141
$ws.Range("H141").Value = 1411
$ws.Range("I141").Value = 1411
$ws.Range("K141").Value = 4233
$ws.Range("M141").Value = 947

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 4211.4
$ws.Range("I9").Value = 320.66666
$ws.Range("J9").Value = 10047.5
$ws.Range("K9").Value = 320.66666
$ws.Range("L9").Value = 10047.5
$ws.Range("M9").Value = -150.66666
$ws.Range("N9").Value = -10387.5
# Row 122
$ws.Range("H122").Value = 6946629.5
$ws.Range("I122").Value = 8334871.5
$ws.Range("K122").Value = 25004614.5
$ws.Range("M122").Value = -25002164.5
# Row 132
$ws.Range("H132").Value = 497
$ws.Range("I132").Value = 497
$ws.Range("K132").Value = 1491
$ws.Range("M132").Value = 1039

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 674
$ws.Range("I9").Value = 674
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 674
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -450
# Row 16
$ws.Range("H16").Value = 962.75
$ws.Range("I16").Value = 351
$ws.Range("K16").Value = 351
$ws.Range("M16").Value = -181
# Row 22
$ws.Range("H22").Value = 2142.2856
$ws.Range("I22").Value = 2142.2856
$ws.Range("K22").Value = 2142.2856
$ws.Range("M22").Value = -1847.2856
# Row 27
$ws.Range("H27").Value = 2142.2856
$ws.Range("I27").Value = 2142.2856
$ws.Range("K27").Value = 2142.2856
$ws.Range("M27").Value = -2035.2856
# Row 30
$ws.Range("H30").Value = 735
$ws.Range("I30").Value = 735
$ws.Range("K30").Value = 735
$ws.Range("M30").Value = -627
# Row 46
$ws.Range("H46").Value = 2583.1292
$ws.Range("I46").Value = 1823.5
$ws.Range("J46").Value = 2847.348
$ws.Range("K46").Value = 1823.5
$ws.Range("L46").Value = 2847.348
$ws.Range("M46").Value = -1635.5
$ws.Range("N46").Value = -3223.348
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
# Row 132
$ws.Range("H132").Value = 8300
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 8300
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").Value = 24900
$ws.Range("N132").Value = -29960
# Row 136
$ws.Range("H136").Value = 3501.5
$ws.Range("I136").Value = 3501.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10504.5
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -7954.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1773.1111
$ws.Range("I122").Value = 1619.75
$ws.Range("K122").Value = 4859.25
$ws.Range("M122").Value = -2409.25
# Row 132
$ws.Range("H132").Value = 2044.2778
$ws.Range("I132").Value = 1982.2858
$ws.Range("K132").Value = 5946.857400000001
$ws.Range("M132").Value = -3416.857400000001
